$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J3").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null
$ws.Range("K3").Value = 2023
$ws.Range("K3").Borders.Item(10).LineStyle = 1

$ws.Range("J4").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null
$ws.Range("K4").Value = 441
$ws.Range("K4").Borders.Item(10).LineStyle = 1

$ws.Range("J5").Copy() | Out-Null
$ws.Range("K5").PasteSpecial(-4122) | Out-Null
$ws.Range("K5").Value = 97
$ws.Range("K5").Borders.Item(10).LineStyle = 1

$ws.Range("J6").Copy() | Out-Null
$ws.Range("K6").PasteSpecial(-4122) | Out-Null
$ws.Range("K6").Value = 344
$ws.Range("K6").Borders.Item(10).LineStyle = 1
